$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Footer/layout date placeholder: 09/12/2015 -> 17/01/2016 ---
$layout = $s.CustomLayout
for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $lsh = $layout.Shapes.Item($i)
    if ($lsh.HasTextFrame -eq -1) {
        $ltr = $lsh.TextFrame.TextRange
        if ($ltr.Text -eq "09/12/2015") {
            $ltr.Text = "17/01/2016"
        }
    }
}

# --- helper: find a shape on slide 1 by its exact current text ---
function Get-ShapeByText($slide, $targetText) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            if ($sh.TextFrame.TextRange.Text -eq $targetText) {
                return $sh
            }
        }
    }
    return $null
}

# --- 2) "nodo de cuarto nivel" -> "ampliaciones y reducciones" ---
$nodoShape = Get-ShapeByText $s "nodo de cuarto nivel"
if ($nodoShape -ne $null) {
    $nodoShape.TextFrame.TextRange.Text = "ampliaciones y reducciones"
}

# --- 3) "k > 1 / k < 1 / k = 1" box: italicize the "k" in each paragraph ---
$kBoxShape = $s.Shapes.Item(27)
$kTr = $kBoxShape.TextFrame.TextRange
if ($kTr.Paragraphs(1).Text.IndexOf("ampliaci") -lt 0) {
    $kBoxShape = $null
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame -eq -1 -and $sh.TextFrame.TextRange.Text.IndexOf("ampliaci") -ge 0 -and $sh.TextFrame.TextRange.Text.IndexOf("congruencia") -ge 0) {
            $kBoxShape = $sh
        }
    }
    $kTr = $kBoxShape.TextFrame.TextRange
}

function Italicize-KInParagraph($shapeTextRange, $paraIndex) {
    $para = $shapeTextRange.Paragraphs($paraIndex)
    $t = $para.Text
    $idx = $t.IndexOf("si k")
    if ($idx -ge 0) {
        $kPosGlobal = $para.Start + $idx + 3
        $kChar = $shapeTextRange.Characters($kPosGlobal, 1)
        $kChar.Font.Italic = -1
    }
}

for ($pi = 1; $pi -le $kTr.Paragraphs().Count; $pi++) {
    Italicize-KInParagraph $kTr $pi
}

# --- 4) "areas a razon k2" box: italicize the superscript "2" ---
$k2Shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -eq -1) {
        $tx = $sh.TextFrame.TextRange.Text
        if ($tx.IndexOf("razón k") -ge 0 -or $tx.IndexOf("razon k") -ge 0) {
            $k2Shape = $sh
        }
    }
}
if ($k2Shape -ne $null) {
    $k2Tr = $k2Shape.TextFrame.TextRange
    $full = $k2Tr.Text
    $lastChar = $k2Tr.Characters($full.Length, 1)
    $lastChar.Font.Italic = -1
}

Write-Host "edit.ps1 completed"
